$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 421 (shifts existing rows 421-469 down to 422-470)
$ws.Rows.Item(421).Insert()

# Populate the newly inserted row 421 with the new record
$ws.Cells.Item(421, 1).Value = 7
$ws.Cells.Item(421, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(421, 3).Value = "Ñuble"
$ws.Cells.Item(421, 4).Value = 45013
$ws.Cells.Item(421, 5).Value = 16
$ws.Cells.Item(421, 6).Value = "Fruta"
$ws.Cells.Item(421, 7).Value = 100103
$ws.Cells.Item(421, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(421, 9).Value = 100103006
$ws.Cells.Item(421, 10).Value = "Nectarín"
$ws.Cells.Item(421, 11).Value = "June Pearl"
$ws.Cells.Item(421, 12).Value = "Primera"
$ws.Cells.Item(421, 13).Value = 120
$ws.Cells.Item(421, 14).Value = 17000
$ws.Cells.Item(421, 15).Value = 18000
$ws.Cells.Item(421, 16).Value = 17500
$ws.Cells.Item(421, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(421, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(421, 19).Value = 1094
$ws.Cells.Item(421, 20).Value = 16
